$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("L6").Value = 1.5
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 2
$ws.Range("N9").Value = 0
